$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (it currently sits after
#    the Picture 20 drawing run; the edit relocates it into the new
#    header block added below).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Mark the two drawing runs (Picture 13 and Picture 20) that are
#    currently missing <w:noProof/> in their run properties.
# ------------------------------------------------------------------
$d.InlineShapes.Item(7).Range.NoProofing = $true
$d.InlineShapes.Item(8).Range.NoProofing = $true

# ------------------------------------------------------------------
# 3. Insert four new paragraphs at the very top of the document: the
#    author name, roll number (carrying the relocated "_GoBack"
#    bookmark), date, and a trailing blank (bold+underlined) line.
# ------------------------------------------------------------------
$first = $d.Paragraphs(1).Range
$first.InsertParagraphBefore()
$top = $d.Paragraphs(1).Range

$headerXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="right"/><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>G P ABINAYA</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="right"/><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>192211241</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="right"/><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>22/02/2023</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
'@

$top.InsertXML($headerXml)
